$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 152383.33
$ws.Range("I52").Value = 152383.33
$ws.Range("K52").Value = 457149.99
$ws.Range("M52").Value = -456989.99
$ws.Range("H98").Value = 157155.4
$ws.Range("I98").Value = 1000
$ws.Range("J98").Value = 261259
$ws.Range("K98").Value = 1000
$ws.Range("L98").Value = 261259
$ws.Range("M98").Value = 498
$ws.Range("N98").Value = -264255
$ws.Range("H122").Value = 157155.4
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 261259
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 783777
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -788677
$ws.Range("H137").Value = 6893.5127
$ws.Range("I137").Value = 7332.3687
$ws.Range("J137").Value = 6476.6
$ws.Range("K137").Value = 21997.1061
$ws.Range("L137").Value = 19429.8
$ws.Range("M137").Value = -19447.1061
$ws.Range("N137").Value = -24529.8
$ws.Range("H141").Value = 2765.2
$ws.Range("I141").Value = 1507.8334
$ws.Range("J141").Value = 4651.25
$ws.Range("K141").Value = 4523.5002
$ws.Range("L141").Value = 13953.75
$ws.Range("M141").Value = 656.4997999999996
$ws.Range("N141").Value = -24313.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10958.211
$ws.Range("I32").Value = 10727.866
$ws.Range("J32").Value = 11822
$ws.Range("K32").Value = 10727.866
$ws.Range("L32").Value = 11822
$ws.Range("M32").Value = -10440.866
$ws.Range("N32").Value = -12396
$ws.Range("H61").Value = 2562.2
$ws.Range("I61").Value = 2030.5454
$ws.Range("J61").Value = 3070.739
$ws.Range("K61").Value = 2030.5454
$ws.Range("L61").Value = 3070.739
$ws.Range("M61").Value = -1818.5454
$ws.Range("N61").Value = -3494.739
$ws.Range("H74").Value = 1584.1621
$ws.Range("I74").Value = 1249.6207
$ws.Range("K74").Value = 1249.6207
$ws.Range("M74").Value = -375.6206999999999
$ws.Range("H77").Value = 1584.1621
$ws.Range("I77").Value = 1249.6207
$ws.Range("K77").Value = 6248.103499999999
$ws.Range("M77").Value = -1880.103499999999
$ws.Range("H132").Value = 3201.0527
$ws.Range("I132").Value = 2340.2693
$ws.Range("J132").Value = 5066.0835
$ws.Range("K132").Value = 7020.8079
$ws.Range("L132").Value = 15198.2505
$ws.Range("M132").Value = -4490.8079
$ws.Range("N132").Value = -20258.2505
$ws.Range("H136").Value = 2562.2
$ws.Range("I136").Value = 2030.5454
$ws.Range("J136").Value = 3070.739
$ws.Range("K136").Value = 6091.6362
$ws.Range("L136").Value = 9212.217000000001
$ws.Range("M136").Value = -3541.6362
$ws.Range("N136").Value = -14312.217

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H105").Value = 2977.2144
$ws.Range("I105").Value = 2525
$ws.Range("J105").Value = 3580.1667
$ws.Range("K105").Value = 2525
$ws.Range("L105").Value = 3580.1667
$ws.Range("M105").Value = -778
$ws.Range("N105").Value = -7074.1667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2017.8334
$ws.Range("I99").Value = 1722.4
$ws.Range("J99").Value = 2228.8572
$ws.Range("K99").Value = 1722.4
$ws.Range("L99").Value = 2228.8572
$ws.Range("M99").Value = -224.4000000000001
$ws.Range("N99").Value = -5224.8572
$ws.Range("H121").Value = 40385.75
$ws.Range("J121").Value = 40385.75
$ws.Range("L121").Value = 40385.75
$ws.Range("N121").Value = -43005.75
$ws.Range("H122").Value = 241072.8
$ws.Range("I122").Value = 400850
$ws.Range("J122").Value = 1407
$ws.Range("K122").Value = 1202550
$ws.Range("L122").Value = 4221
$ws.Range("M122").Value = -1200100
$ws.Range("N122").Value = -9121
$ws.Range("H126").Value = 2017.8334
$ws.Range("I126").Value = 1722.4
$ws.Range("J126").Value = 2228.8572
$ws.Range("K126").Value = 5167.200000000001
$ws.Range("L126").Value = 6686.571599999999
$ws.Range("M126").Value = -2697.200000000001
$ws.Range("N126").Value = -11626.5716
$ws.Range("H134").Value = 484122.06
$ws.Range("I134").Value = 1137.72
$ws.Range("J134").Value = 3502774.2
$ws.Range("K134").Value = 3413.16
$ws.Range("L134").Value = 10508322.6
$ws.Range("M134").Value = -878.1599999999999
$ws.Range("N134").Value = -10513392.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 23368.385
$ws.Range("I18").Value = 25290.75
$ws.Range("K18").Value = 75872.25
$ws.Range("M18").Value = -75703.25
$ws.Range("H112").Value = 2687.5
$ws.Range("I112").Value = 326.66666
$ws.Range("J112").Value = 4104
$ws.Range("K112").Value = 979.9999799999999
$ws.Range("L112").Value = 12312
$ws.Range("M112").Value = 128.0000200000001
$ws.Range("N112").Value = -14528

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1925
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 30496.715
$ws.Range("I126").Value = 51651.5
$ws.Range("K126").Value = 154954.5
$ws.Range("M126").Value = -152484.5
$ws.Range("H132").Value = 3346
$ws.Range("I132").Value = 2405.5
$ws.Range("J132").Value = 5227
$ws.Range("K132").Value = 7216.5
$ws.Range("L132").Value = 15681
$ws.Range("M132").Value = -4686.5
$ws.Range("N132").Value = -20741

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 66670400
$ws.Range("I7").Value = 100002480
$ws.Range("J7").Value = 6241
$ws.Range("K7").Value = 100002480
$ws.Range("L7").Value = 6241
$ws.Range("M7").Value = -100002368
$ws.Range("N7").Value = -6465
$ws.Range("H16").Value = 2697.75
$ws.Range("I16").Value = 2843
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 2843
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -2673
$ws.Range("N16").Value = -1440
$ws.Range("H40").Value = 6510.7
$ws.Range("I40").Value = 7101
$ws.Range("K40").Value = 7101
$ws.Range("M40").Value = -6965
$ws.Range("H126").Value = 66670400
$ws.Range("I126").Value = 100002480
$ws.Range("J126").Value = 6241
$ws.Range("K126").Value = 300007440
$ws.Range("L126").Value = 18723
$ws.Range("M126").Value = -300004970
$ws.Range("N126").Value = -23663
